$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = '@'
$r.Value = '62.176.81'
$r.Style = 'Normal'
$r = $ws.Range('E2')
$r.NumberFormat = '@'
$r.Value = '  -2.56%  '
$r.Style = 'Normal'
$r = $ws.Range('D3')
$r.NumberFormat = '@'
$r.Value = '2.995.63'
$r.Style = 'Normal'
$r = $ws.Range('E3')
$r.NumberFormat = '@'
$r.Value = '  -2.85%  '
$r.Style = 'Normal'
$r = $ws.Range('E4')
$r.NumberFormat = '@'
$r.Value = '  +0.09%  '
$r.Style = 'Normal'
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '580.34'
$r.Style = 'Normal'
$r = $ws.Range('E5')
$r.NumberFormat = '@'
$r.Value = '  -1.51%  '
$r.Style = 'Normal'
$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '146.30'
$r.Style = 'Normal'
$r = $ws.Range('E6')
$r.NumberFormat = '@'
$r.Value = '  -5.99%  '
$r.Style = 'Normal'
$r = $ws.Range('E7')
$r.NumberFormat = '@'
$r.Value = '  +0.06%  '
$r.Style = 'Normal'
$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '0.521'
$r.Style = 'Normal'
$r = $ws.Range('E8')
$r.NumberFormat = '@'
$r.Value = '  -3.65%  '
$r.Style = 'Normal'
$r = $ws.Range('D9')
$r.NumberFormat = '@'
$r.Value = '2.992.10'
$r.Style = 'Normal'
$r = $ws.Range('E9')
$r.NumberFormat = '@'
$r.Value = '  -2.98%  '
$r.Style = 'Normal'
$r = $ws.Range('E10')
$r.NumberFormat = '@'
$r.Value = '  -5.60%  '
$r.Style = 'Normal'
$r = $ws.Range('E11')
$r.NumberFormat = '@'
$r.Value = '  -3.95%  '
$r.Style = 'Normal'
$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '0.440'
$r.Style = 'Normal'
$r = $ws.Range('E12')
$r.NumberFormat = '@'
$r.Value = '  -2.35%  '
$r.Style = 'Normal'
$r = $ws.Range('E13')
$r.NumberFormat = '@'
$r.Value = '  -4.67%  '
$r.Style = 'Normal'
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '34.45'
$r.Style = 'Normal'
$r = $ws.Range('E14')
$r.NumberFormat = '@'
$r.Value = '  -5.97%  '
$r.Style = 'Normal'
$r = $ws.Range('E15')
$r.NumberFormat = '@'
$r.Value = '  +1.59%  '
$r.Style = 'Normal'
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '3.487.59'
$r.Style = 'Normal'
$r = $ws.Range('E16')
$r.NumberFormat = '@'
$r.Value = '  -2.80%  '
$r.Style = 'Normal'
$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '7.02'
$r.Style = 'Normal'
$r = $ws.Range('E17')
$r.NumberFormat = '@'
$r.Value = '  -2.16%  '
$r.Style = 'Normal'
$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '62.218.21'
$r.Style = 'Normal'
$r = $ws.Range('E18')
$r.NumberFormat = '@'
$r.Value = '  -2.22%  '
$r.Style = 'Normal'
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '2.994.88'
$r.Style = 'Normal'
$r = $ws.Range('E19')
$r.NumberFormat = '@'
$r.Value = '  -2.78%  '
$r.Style = 'Normal'
$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '451.60'
$r.Style = 'Normal'
$r = $ws.Range('E20')
$r.NumberFormat = '@'
$r.Value = '  -3.92%  '
$r.Style = 'Normal'
$r = $ws.Range('E21')
$r.NumberFormat = '@'
$r.Value = '  -3.53%  '
$r.Style = 'Normal'
$r = $ws.Range('E22')
$r.NumberFormat = '@'
$r.Value = '  -4.03%  '
$r.Style = 'Normal'
$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '7.28'
$r.Style = 'Normal'
$r = $ws.Range('E23')
$r.NumberFormat = '@'
$r.Value = '  -2.86%  '
$r.Style = 'Normal'
$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '2.29'
$r.Style = 'Normal'
$r = $ws.Range('E24')
$r.NumberFormat = '@'
$r.Value = '  -5.44%  '
$r.Style = 'Normal'
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '79.95'
$r.Style = 'Normal'
$r = $ws.Range('E25')
$r.NumberFormat = '@'
$r.Value = '  -0.65%  '
$r.Style = 'Normal'
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '12.22'
$r.Style = 'Normal'
$r = $ws.Range('E26')
$r.NumberFormat = '@'
$r.Value = '  -4.88%  '
$r.Style = 'Normal'
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '10.01'
$r.Style = 'Normal'
$r = $ws.Range('E27')
$r.NumberFormat = '@'
$r.Value = '  -3.44%  '
$r.Style = 'Normal'
$r = $ws.Range('E28')
$r.NumberFormat = '@'
$r.Value = '  -0.04%  '
$r.Style = 'Normal'
$r = $ws.Range('E29')
$r.NumberFormat = '@'
$r.Value = '  +0.16%  '
$r.Style = 'Normal'
$r = $ws.Range('B30')
$r.NumberFormat = '@'
$r.Value = 'NEARProtocol'
$r.Style = 'Normal'
$r = $ws.Range('C30')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$r.Style = 'Normal'
$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '7.13'
$r.Style = 'Normal'
$r = $ws.Range('E30')
$r.NumberFormat = '@'
$r.Value = '  -3.41%  '
$r.Style = 'Normal'
$r = $ws.Range('B31')
$r.NumberFormat = '@'
$r.Value = 'PancakeSwap'
$r.Style = 'Normal'
$r = $ws.Range('C31')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$r.Style = 'Normal'
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '2.60'
$r.Style = 'Normal'
$r = $ws.Range('E31')
$r.NumberFormat = '@'
$r.Value = '  -2.22%  '
$r.Style = 'Normal'
$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '2.08'
$r.Style = 'Normal'
$r = $ws.Range('E32')
$r.NumberFormat = '@'
$r.Value = '  -3.11%  '
$r.Style = 'Normal'
$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '26.73'
$r.Style = 'Normal'
$r = $ws.Range('E33')
$r.NumberFormat = '@'
$r.Value = '  -1.42%  '
$r.Style = 'Normal'
$r = $ws.Range('E34')
$r.NumberFormat = '@'
$r.Value = '  -5.44%  '
$r.Style = 'Normal'
$r = $ws.Range('E35')
$r.NumberFormat = '@'
$r.Value = '  -2.85%  '
$r.Style = 'Normal'
$r = $ws.Range('D36')
$r.NumberFormat = '@'
$r.Value = '0.0₃0787'
$r.Style = 'Normal'
$r = $ws.Range('E36')
$r.NumberFormat = '@'
$r.Value = '  -5.49%  '
$r.Style = 'Normal'
$r = $ws.Range('E37')
$r.NumberFormat = '@'
$r.Value = '  -4.35%  '
$r.Style = 'Normal'
$r = $ws.Range('E38')
$r.NumberFormat = '@'
$r.Value = '  -4.58%  '
$r.Style = 'Normal'
$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '50.08'
$r.Style = 'Normal'
$r = $ws.Range('E39')
$r.NumberFormat = '@'
$r.Value = '  -1.07%  '
$r.Style = 'Normal'
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '8.98'
$r.Style = 'Normal'
$r = $ws.Range('E40')
$r.NumberFormat = '@'
$r.Value = '  -1.41%  '
$r.Style = 'Normal'
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '2.93'
$r.Style = 'Normal'
$r = $ws.Range('E41')
$r.NumberFormat = '@'
$r.Value = '  -9.74%  '
$r.Style = 'Normal'
$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '407.45'
$r.Style = 'Normal'
$r = $ws.Range('E42')
$r.NumberFormat = '@'
$r.Value = '  -6.36%  '
$r.Style = 'Normal'
$r = $ws.Range('B43')
$r.NumberFormat = '@'
$r.Value = 'TheGraph'
$r.Style = 'Normal'
$r = $ws.Range('C43')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$r.Style = 'Normal'
$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '0.276'
$r.Style = 'Normal'
$r = $ws.Range('E43')
$r.NumberFormat = '@'
$r.Value = '  -4.49%  '
$r.Style = 'Normal'
$r = $ws.Range('B44')
$r.NumberFormat = '@'
$r.Value = 'Kaspa'
$r.Style = 'Normal'
$r = $ws.Range('C44')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$r.Style = 'Normal'
$r = $ws.Range('D44')
$r.NumberFormat = '@'
$r.Value = '0.111'
$r.Style = 'Normal'
$r = $ws.Range('E44')
$r.NumberFormat = '@'
$r.Value = '  -0.57%  '
$r.Style = 'Normal'
$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '2.758.45'
$r.Style = 'Normal'
$r = $ws.Range('E45')
$r.NumberFormat = '@'
$r.Value = '  -2.07%  '
$r.Style = 'Normal'
$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '0.0350'
$r.Style = 'Normal'
$r = $ws.Range('E46')
$r.NumberFormat = '@'
$r.Value = '  -2.87%  '
$r.Style = 'Normal'
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '38.14'
$r.Style = 'Normal'
$r = $ws.Range('E47')
$r.NumberFormat = '@'
$r.Value = '  -4.22%  '
$r.Style = 'Normal'
$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '127.59'
$r.Style = 'Normal'
$r = $ws.Range('E48')
$r.NumberFormat = '@'
$r.Value = '  -1.74%  '
$r.Style = 'Normal'
$r = $ws.Range('E49')
$r.NumberFormat = '@'
$r.Value = '  +0.04%  '
$r.Style = 'Normal'
$r = $ws.Range('E50')
$r.NumberFormat = '@'
$r.Value = '  -1.91%  '
$r.Style = 'Normal'
$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '23.71'
$r.Style = 'Normal'
$r = $ws.Range('E51')
$r.NumberFormat = '@'
$r.Value = '  -4.74%  '
$r.Style = 'Normal'
